$wb = $excel.ActiveWorkbook

foreach ($name in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a new column before the "T.Tremble" column (R) and populate the
    # new player's header/data, shifting the remaining headers one column right.
    $ws.Range("R1").EntireColumn.Insert()
    $ws.Range("R1").Value = "C.Saunders"
    $ws.Range("R2").Value = "n"
}
